# Generate Report for Handoff
# Updates the Overview / zh-cn / de-de sheets with the new handoff batch:
#   - "66f80cc9-...md" -> "dd0f5e15-...md" (status now "Ready for handoff")
#   - new source file "ffff4158f09c-...md" (status "Ready for handoff")
#   - new handoff (.xlf) artifacts recorded for zh-cn / de-de
#   - ".localization-config" row moves down to row 4

$wb = $excel.ActiveWorkbook

$oldMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f95aa65571649722770d9a6ad296f3a68a164924/e2e"
$newMdName1 = "dd0f5e15-c13b-4852-8032-8f74ba7423cc.md"
$newMdName2 = "ffff4158f09c-d6ae-4da2-b688-eff234e29588.md"
$configName = ".localization-config"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/36d61a40b8ec7f1339fc85ecbc903c4c7b672d31/.localization-config"

$xlfZhCn = "dd0f5e15-c13b-4852-8032-8f74ba7423cc.1f5f16f89e92bfdad36837929344c65010bf9e84.zh-cn.xlf"
$xlfDeDe = "dd0f5e15-c13b-4852-8032-8f74ba7423cc.1f5f16f89e92bfdad36837929344c65010bf9e84.de-de.xlf"
$xlfUrlBase = "https://github.com/OpenLocalizationTest/ol-handback/blob/1f5f16f89e92bfdad36837929344c65010bf9e84"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMdName1
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

$ws.Range("A3").Value = $newMdName2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = $configName
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), "$oldMdUrl/$newMdName1", "", "", $newMdName1)
$ws.Hyperlinks.Add($ws.Range("A3"), "$oldMdUrl/$newMdName2", "", "", $newMdName2)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $configName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMdName1
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $xlfZhCn
$ws.Range("D2").Value = "2016-02-17 04:46:46"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = $newMdName2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $xlfZhCn
$ws.Range("D3").Value = "2016-02-17 04:46:46"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = $configName
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$oldMdUrl/$newMdName1", "", "", $newMdName1)
$ws.Hyperlinks.Add($ws.Range("C2"), "$xlfUrlBase/$xlfZhCn", "", "", $xlfZhCn)
$ws.Hyperlinks.Add($ws.Range("A3"), "$oldMdUrl/$newMdName2", "", "", $newMdName2)
$ws.Hyperlinks.Add($ws.Range("C3"), "$xlfUrlBase/$xlfZhCn", "", "", $xlfZhCn)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $configName)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMdName1
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $xlfDeDe
$ws.Range("D2").Value = "2016-02-17 04:46:57"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = $newMdName2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $xlfDeDe
$ws.Range("D3").Value = "2016-02-17 04:46:57"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = $configName
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$oldMdUrl/$newMdName1", "", "", $newMdName1)
$ws.Hyperlinks.Add($ws.Range("C2"), "$xlfUrlBase/$xlfDeDe", "", "", $xlfDeDe)
$ws.Hyperlinks.Add($ws.Range("A3"), "$oldMdUrl/$newMdName2", "", "", $newMdName2)
$ws.Hyperlinks.Add($ws.Range("C3"), "$xlfUrlBase/$xlfDeDe", "", "", $xlfDeDe)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $configName)
